$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = New-Object 'object[,]' 24,1
$colB[0,0] = 0.9744635030973541
$colB[1,0] = 0.9073526923759516
$colB[2,0] = 0.8663852883026948
$colB[3,0] = 0.84975164757887
$colB[4,0] = 0.8469933477251743
$colB[5,0] = 0.8661607134830831
$colB[6,0] = 0.9512745869183732
$colB[7,0] = 1.120050785307001
$colB[8,0] = 1.245166245880966
$colB[9,0] = 1.302322796915007
$colB[10,0] = 1.324000521203118
$colB[11,0] = 1.319330346708512
$colB[12,0] = 1.304105564263239
$colB[13,0] = 1.294784318995141
$colB[14,0] = 1.241435704438118
$colB[15,0] = 1.208769140981303
$colB[16,0] = 1.190002923285476
$colB[17,0] = 1.183652934005011
$colB[18,0] = 1.212244206691992
$colB[19,0] = 1.308576541994285
$colB[20,0] = 1.371731788082172
$colB[21,0] = 1.338006947520228
$colB[22,0] = 1.210673085111125
$colB[23,0] = 1.074194797741995
$ws.Range("B2:B25").Value = $colB

$colD = New-Object 'object[,]' 24,1
$colD[0,0] = 0.004378622681091571
$colD[1,0] = 0.004253294322388967
$colD[2,0] = 0.004183880374812432
$colD[3,0] = 0.004157450779707972
$colD[4,0] = 0.00415317321596298
$colD[5,0] = 0.004183516466351733
$colD[6,0] = 0.004333821685982286
$colD[7,0] = 0.004690138858649107
$colD[8,0] = 0.004991929501759529
$colD[9,0] = 0.005138449496261188
$colD[10,0] = 0.005195304905068809
$colD[11,0] = 0.005182998378231218
$colD[12,0] = 0.005143099286691211
$colD[13,0] = 0.005118839850808854
$colD[14,0] = 0.004982543603123446
$colD[15,0] = 0.004901326178693211
$colD[16,0] = 0.004855478562461712
$colD[17,0] = 0.004840102848753247
$colD[18,0] = 0.004909881918983672
$colD[19,0] = 0.005154781045025203
$colD[20,0] = 0.005322851382768334
$colD[21,0] = 0.005232401181721968
$colD[22,0] = 0.004906011239405927
$colD[23,0] = 0.00458690071072354
$ws.Range("D2:D25").Value = $colD

$colE = New-Object 'object[,]' 24,1
$colE[0,0] = 0.4483793404404963
$colE[1,0] = 0.4222477322751601
$colE[2,0] = 0.406099570183585
$colE[3,0] = 0.3994930833390313
$colE[4,0] = 0.3983945117879841
$colE[5,0] = 0.4060105779893703
$colE[6,0] = 0.4393905719355971
$colE[7,0] = 0.5040334991671216
$colE[8,0] = 0.5510402301691499
$colE[9,0] = 0.5723218452618397
$colE[10,0] = 0.5803660988974997
$colE[11,0] = 0.5786342751075892
$colE[12,0] = 0.5729839438815532
$colE[13,0] = 0.5695210484367692
$colE[14,0] = 0.549647377046
$colE[15,0] = 0.5374294239915685
$colE[16,0] = 0.5303923789913654
$colE[17,0] = 0.5280081065289437
$colE[18,0] = 0.5387310394632863
$colE[19,0] = 0.5746439799387986
$colE[20,0] = 0.598029896489237
$colE[21,0] = 0.5855561815750576
$colE[22,0] = 0.5381426189206877
$colE[23,0] = 0.4866320816636147
$ws.Range("E2:E25").Value = $colE

$colF = New-Object 'object[,]' 24,1
$colF[0,0] = 0.9127886060284141
$colF[1,0] = 0.8536296548066389
$colF[2,0] = 0.8180078375319937
$colF[3,0] = 0.8036660063232972
$colF[4,0] = 0.801295023036019
$colF[5,0] = 0.8178137155276062
$colF[6,0] = 0.8922435600163254
$colF[7,0] = 1.043877345674346
$colF[8,0] = 1.158903968361727
$colF[9,0] = 1.21205451038665
$colF[10,0] = 1.23230245313178
$colF[11,0] = 1.227936276252535
$colF[12,0] = 1.213717881975356
$colF[13,0] = 1.205024533936381
$colF[14,0] = 1.155447251894856
$colF[15,0] = 1.125246062582562
$colF[16,0] = 1.107952581513871
$colF[17,0] = 1.102110533059317
$colF[18,0] = 1.128452997941736
$colF[19,0] = 1.217890865879525
$colF[20,0] = 1.27705003724526
$colF[21,0] = 1.245410227232981
$colF[22,0] = 1.127002925908641
$colF[23,0] = 1.002232136738485
$ws.Range("F2:F25").Value = $colF

$colG = New-Object 'object[,]' 24,1
$colG[0,0] = 0.8329604128828976
$colG[1,0] = 0.7651171033127184
$colG[2,0] = 0.7240857060200199
$colG[3,0] = 0.7075195461595172
$colG[4,0] = 0.7047780014643195
$colG[5,0] = 0.7238616669987152
$colG[6,0] = 0.8094367925291124
$colG[7,0] = 0.9823328087415177
$colG[8,0] = 1.112647672837909
$colG[9,0] = 1.172687172392443
$colG[10,0] = 1.195534819822853
$colG[11,0] = 1.190609142692921
$colG[12,0] = 1.174564600861288
$colG[13,0] = 1.164751527115129
$colG[14,0] = 1.108739463132565
$colG[15,0] = 1.074574231782663
$colG[16,0] = 1.054994523156751
$colG[17,0] = 1.048377315886938
$colG[18,0] = 1.078203786138175
$colG[19,0] = 1.179274208988659
$colG[20,0] = 1.245983718765672
$colG[21,0] = 1.210318775586614
$colG[22,0] = 1.076562670664259
$colG[23,0] = 0.9349961907399518
$ws.Range("G2:G25").Value = $colG

$colH = New-Object 'object[,]' 24,1
$colH[0,0] = 0.7074198953251312
$colH[1,0] = 0.6817852656444643
$colH[2,0] = 0.6665915257456732
$colH[3,0] = 0.6605361874273399
$colH[4,0] = 0.6595389002716558
$colH[5,0] = 0.6665093111498663
$colH[6,0] = 0.6984670916660889
$colH[7,0] = 0.7655212618909104
$colH[8,0] = 0.8175383477593527
$colH[9,0] = 0.8418175546848943
$colH[10,0] = 0.8511014210631345
$colH[11,0] = 0.8490979582183797
$colH[12,0] = 0.8425795360104189
$colH[13,0] = 0.8385985529104687
$colH[14,0] = 0.8159641517485738
$colH[15,0] = 0.8022374135012171
$colH[16,0] = 0.7944000985766593
$colH[17,0] = 0.7917564329430604
$colH[18,0] = 0.8036926426473485
$colH[19,0] = 0.8444917060750754
$colH[20,0] = 0.8716805913802546
$colH[21,0] = 0.8571209858846203
$colH[22,0] = 0.8030345642927728
$colH[23,0] = 0.7469037561691607
$ws.Range("H2:H25").Value = $colH

$colI = New-Object 'object[,]' 24,1
$colI[0,0] = 0.6728402859271183
$colI[1,0] = 0.6863084089886957
$colI[2,0] = 0.6950605011788689
$colI[3,0] = 0.6987481819951631
$colI[4,0] = 0.6993678281274223
$colI[5,0] = 0.6951097443734504
$colI[6,0] = 0.6773837956157704
$colI[7,0] = 0.6464620844305626
$colI[8,0] = 0.626096502148286
$colI[9,0] = 0.6173451677675015
$colI[10,0] = 0.6141052908263411
$colI[11,0] = 0.6147997577506406
$colI[12,0] = 0.6170771351636084
$colI[13,0] = 0.6184817483106873
$colI[14,0] = 0.6266787710055581
$colI[15,0] = 0.6318389905210449
$colI[16,0] = 0.6348552878772615
$colI[17,0] = 0.635884837069236
$colI[18,0] = 0.6312846779532535
$colI[19,0] = 0.6164062018750034
$colI[20,0] = 0.6071140895847495
$colI[21,0] = 0.612033857635236
$colI[22,0] = 0.6315351281363935
$colI[23,0] = 0.6544149851221244
$ws.Range("I2:I25").Value = $colI

$colL = New-Object 'object[,]' 24,1
$colL[0,0] = 0.4009955759892705
$colL[1,0] = 0.3591094937132766
$colL[2,0] = 0.333474467128184
$colL[3,0] = 0.323048803636965
$colL[4,0] = 0.3213188847081483
$colL[5,0] = 0.3333337788231461
$colL[6,0] = 0.3865359103829462
$colL[7,0] = 0.4915356429366113
$colL[8,0] = 0.5691103904212582
$colL[9,0] = 0.6045003337608534
$colL[10,0] = 0.617916359154691
$colL[11,0] = 0.6150263256179187
$colL[12,0] = 0.6056037844423088
$colL[13,0] = 0.5998341140168293
$colL[14,0] = 0.5667996134211819
$colL[15,0] = 0.5465599224939695
$colL[16,0] = 0.5349280674062129
$colL[17,0] = 0.530991345583061
$colL[18,0] = 0.5487134901289608
$colL[19,0] = 0.6083710153037316
$colL[20,0] = 0.6474459494284304
$colL[21,0] = 0.6265830732639301
$colL[22,0] = 0.5477398490885719
$colL[23,0] = 0.4630563351569776
$ws.Range("L2:L25").Value = $colL

$colN = New-Object 'object[,]' 24,1
$colN[0,0] = 1.311147546546579
$colN[1,0] = 1.305393740739518
$colN[2,0] = 1.302396190758202
$colN[3,0] = 1.301309910960228
$colN[4,0] = 1.301137722497202
$colN[5,0] = 1.302380992354756
$colN[6,0] = 1.309052905668224
$colN[7,0] = 1.326358302800884
$colN[8,0] = 1.341615317907696
$colN[9,0] = 1.349102077971622
$colN[10,0] = 1.352015116055796
$colN[11,0] = 1.351384282941297
$colN[12,0] = 1.349340176174465
$colN[13,0] = 1.348098237764987
$colN[14,0] = 1.341136982319171
$colN[15,0] = 1.337005968014509
$colN[16,0] = 1.334681369618409
$colN[17,0] = 1.333903155517689
$colN[18,0] = 1.337440400770575
$colN[19,0] = 1.349938468581087
$colN[20,0] = 1.358560849786457
$colN[21,0] = 1.353917557546225
$colN[22,0] = 1.337243836764472
$colN[23,0] = 1.321228005521021
$ws.Range("N2:N25").Value = $colN

